$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed rows 6-11 from row 2 so date/percent number formats (styles) carry over correctly,
# then overwrite each cell with the real values from the new data.
$ws.Range("A2:Y2").Copy($ws.Range("A6:Y6"))
$ws.Range("A2:Y2").Copy($ws.Range("A7:Y7"))
$ws.Range("A2:Y2").Copy($ws.Range("A8:Y8"))
$ws.Range("A2:Y2").Copy($ws.Range("A9:Y9"))
$ws.Range("A2:Y2").Copy($ws.Range("A10:Y10"))
$ws.Range("A2:Y2").Copy($ws.Range("A11:Y11"))

$ws.Range("A6").Value = 42650.338333333333
$ws.Range("B6").Value = 11
$ws.Range("C6").Value = "Buy"
$ws.Range("D6").Value = 40
$ws.Range("E6").Value = 5500
$ws.Range("F6").Value = 893
$ws.Range("G6").Value = 68
$ws.Range("H6").Value = 29
$ws.Range("I6").Value = 84
$ws.Range("J6").Value = 15
$ws.Range("K6").Value = 8678
$ws.Range("L6").Value = 149
$ws.Range("M6").Value = 64
$ws.Range("N6").Value = 54
$ws.Range("O6").Value = 10
$ws.Range("P6").Value = "Noun"
$ws.Range("Q6").Value = 38.48959524716075
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0.1046
$ws.Range("T6").Value = 0.034500000000000003
$ws.Range("U6").Value = 4.82
$ws.Range("V6").Value = 2.2799999999999998
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 0.18999999999999773
$ws.Range("Y6").Value = "Up"
$ws.Range("A7").Value = 42650.339583333334
$ws.Range("B7").Value = -10
$ws.Range("C7").Value = "Down"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2025
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = "Noun"
$ws.Range("Q7").Value = 38.48959524716075
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0.1046
$ws.Range("T7").Value = 0.034500000000000003
$ws.Range("U7").Value = 4.82
$ws.Range("V7").Value = 2.2799999999999998
$ws.Range("W7").Value = 0
$ws.Range("X7").Value = 0.18999999999999773
$ws.Range("Y7").Value = "Up"
$ws.Range("A8").Value = 42650.348738425928
$ws.Range("B8").Value = 11
$ws.Range("C8").Value = "Buy"
$ws.Range("D8").Value = 28
$ws.Range("E8").Value = 22967
$ws.Range("F8").Value = 3627
$ws.Range("G8").Value = 61
$ws.Range("H8").Value = 33
$ws.Range("I8").Value = 84
$ws.Range("J8").Value = 14
$ws.Range("K8").Value = 21122
$ws.Range("L8").Value = 514
$ws.Range("M8").Value = 279
$ws.Range("N8").Value = 80
$ws.Range("O8").Value = 14
$ws.Range("P8").Value = "Noun"
$ws.Range("Q8").Value = 38.48959524716075
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0.1046
$ws.Range("T8").Value = 0.034500000000000003
$ws.Range("U8").Value = 4.82
$ws.Range("V8").Value = 2.2799999999999998
$ws.Range("W8").Value = 0
$ws.Range("X8").Value = 0.18999999999999773
$ws.Range("Y8").Value = "Up"
$ws.Range("A9").Value = 42650.359016203707
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "Down"
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 1929
$ws.Range("F9").Value = 315
$ws.Range("G9").Value = 60
$ws.Range("H9").Value = 39
$ws.Range("I9").Value = 50
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 5168
$ws.Range("L9").Value = 41
$ws.Range("M9").Value = 27
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = "Noun"
$ws.Range("Q9").Value = 38.48959524716075
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0.1046
$ws.Range("T9").Value = 0.034500000000000003
$ws.Range("U9").Value = 4.82
$ws.Range("V9").Value = 2.2799999999999998
$ws.Range("W9").Value = 0
$ws.Range("X9").Value = 0.18999999999999773
$ws.Range("Y9").Value = "Up"
$ws.Range("A10").Value = 42650.361435185187
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = "Buy"
$ws.Range("D10").Value = 40
$ws.Range("E10").Value = 5501
$ws.Range("F10").Value = 893
$ws.Range("G10").Value = 68
$ws.Range("H10").Value = 29
$ws.Range("I10").Value = 84
$ws.Range("J10").Value = 15
$ws.Range("K10").Value = 10208
$ws.Range("L10").Value = 149
$ws.Range("M10").Value = 64
$ws.Range("N10").Value = 54
$ws.Range("O10").Value = 10
$ws.Range("P10").Value = "Noun"
$ws.Range("Q10").Value = 38.48959524716075
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0.1046
$ws.Range("T10").Value = 0.034500000000000003
$ws.Range("U10").Value = 4.82
$ws.Range("V10").Value = 2.2799999999999998
$ws.Range("W10").Value = 0
$ws.Range("X10").Value = 0.18999999999999773
$ws.Range("Y10").Value = "Up"
$ws.Range("A11").Value = 42650.363067129627
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "Down"
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 2039
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = "Noun"
$ws.Range("Q11").Value = 37.799019424898844
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0.1046
$ws.Range("T11").Value = 0.034299999999999997
$ws.Range("U11").Value = 4.82
$ws.Range("V11").Value = 2.2799999999999998
$ws.Range("W11").Value = 0
# Row 11 has no PriceChange / UpDown values in the source data - remove the copied placeholders.
$ws.Range("X11:Y11").ClearContents()

# Existing row 5 gains the PriceChange / UpDown columns too.
$ws.Range("X5").Value = 0.18999999999999773
$ws.Range("Y5").Value = "Up"

# Match the author's final selection.
$ws.Range("B7").Select() | Out-Null
